$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell C4 held a numeric price (150) but should instead hold a name
# typed in by mistake - simulate the "wrong type" value arriving from
# the sheet by overwriting it with text.
$ws.Range("C4").Value = "LEANDRO"

# Column C is a "best fit" column (its width auto-adjusts to content);
# re-fit it now that C4 holds a longer text value instead of a short number.
$ws.Columns("C").AutoFit()

# Move the active selection to E6 (matches the post-edit selection).
$ws.Range("E6").Select()
